$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 547, shifting existing rows 547-614 down to 550-617
$ws.Rows("547:549").Insert()

# Row 547
$ws.Range("A547").Value = 2
$ws.Range("B547").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C547").Value = "Coquimbo"
$ws.Range("D547").Value = 45077
$ws.Range("E547").Value = 4
$ws.Range("F547").Value = "Fruta"
$ws.Range("G547").Value = 100101
$ws.Range("H547").Value = "Berries"
$ws.Range("I547").Value = 100112025
$ws.Range("J547").Value = "Frutilla"
$ws.Range("K547").Value = "Sin especificar"
$ws.Range("L547").Value = "Especial"
$ws.Range("M547").Value = 500
$ws.Range("N547").Value = 20000
$ws.Range("O547").Value = 21000
$ws.Range("P547").Value = 20500
$ws.Range("Q547").Value = "`$/bandeja 7 kilos"
$ws.Range("R547").Value = "Provincia de Melipilla"
$ws.Range("S547").Value = 2929
$ws.Range("T547").Value = 7

# Row 548
$ws.Range("A548").Value = 2
$ws.Range("B548").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C548").Value = "Coquimbo"
$ws.Range("D548").Value = 45077
$ws.Range("E548").Value = 4
$ws.Range("F548").Value = "Fruta"
$ws.Range("G548").Value = 100101
$ws.Range("H548").Value = "Berries"
$ws.Range("I548").Value = 100112025
$ws.Range("J548").Value = "Frutilla"
$ws.Range("K548").Value = "Sin especificar"
$ws.Range("L548").Value = "Primera"
$ws.Range("M548").Value = 400
$ws.Range("N548").Value = 17000
$ws.Range("O548").Value = 18000
$ws.Range("P548").Value = 17500
$ws.Range("Q548").Value = "`$/bandeja 7 kilos"
$ws.Range("R548").Value = "Provincia de Melipilla"
$ws.Range("S548").Value = 2500
$ws.Range("T548").Value = 7

# Row 549
$ws.Range("A549").Value = 2
$ws.Range("B549").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C549").Value = "Coquimbo"
$ws.Range("D549").Value = 45077
$ws.Range("E549").Value = 4
$ws.Range("F549").Value = "Fruta"
$ws.Range("G549").Value = 100101
$ws.Range("H549").Value = "Berries"
$ws.Range("I549").Value = 100112025
$ws.Range("J549").Value = "Frutilla"
$ws.Range("K549").Value = "Sin especificar"
$ws.Range("L549").Value = "Segunda"
$ws.Range("M549").Value = 300
$ws.Range("N549").Value = 13000
$ws.Range("O549").Value = 14000
$ws.Range("P549").Value = 13500
$ws.Range("Q549").Value = "`$/bandeja 7 kilos"
$ws.Range("R549").Value = "Provincia de Melipilla"
$ws.Range("S549").Value = 1929
$ws.Range("T549").Value = 7
